# Updates cryptos list figures (price + 1h volume%) to the latest scrape,
# and fixes the row ordering for two coin pairs whose ranking flipped
# (Algorand/EnergySwap at rows 46-47, Cronos/EOS at rows 50-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a target cell + its new text. "Numeric" marks D-column price
# cells whose text would otherwise be auto-parsed into a number by Excel
# (stripping meaningful trailing/leading zeros, e.g. "0.06942" -> 6.942E-2);
# for those we force the Text number format before writing, then restore the
# cell's style so formatting matches the rest of the sheet.
$updates = @(
    @{ Ref = "D2"; Value = "26.480.54"; Numeric = $false },
    @{ Ref = "E2"; Value = "  -0.30%  "; Numeric = $false },
    @{ Ref = "D3"; Value = "1.838.38"; Numeric = $false },
    @{ Ref = "E3"; Value = "  -0.64%  "; Numeric = $false },
    @{ Ref = "E4"; Value = "  -0.05%  "; Numeric = $false },
    @{ Ref = "D5"; Value = "260.48"; Numeric = $true },
    @{ Ref = "E5"; Value = "  -0.97%  "; Numeric = $false },
    @{ Ref = "E6"; Value = "  -0.03%  "; Numeric = $false },
    @{ Ref = "D7"; Value = "0.5378"; Numeric = $true },
    @{ Ref = "E7"; Value = "  +2.44%  "; Numeric = $false },
    @{ Ref = "D8"; Value = "0.2920"; Numeric = $true },
    @{ Ref = "E8"; Value = "  -9.77%  "; Numeric = $false },
    @{ Ref = "D9"; Value = "0.06942"; Numeric = $true },
    @{ Ref = "E9"; Value = "  +1.99%  "; Numeric = $false },
    @{ Ref = "D10"; Value = "17.29"; Numeric = $true },
    @{ Ref = "E10"; Value = "  -8.72%  "; Numeric = $false },
    @{ Ref = "D11"; Value = "1.852.77"; Numeric = $false },
    @{ Ref = "E11"; Value = "  +0.07%  "; Numeric = $false },
    @{ Ref = "D12"; Value = "0.7262"; Numeric = $true },
    @{ Ref = "E12"; Value = "  -7.42%  "; Numeric = $false },
    @{ Ref = "D13"; Value = "0.07195"; Numeric = $true },
    @{ Ref = "E13"; Value = "  -7.37%  "; Numeric = $false },
    @{ Ref = "D14"; Value = "89.01"; Numeric = $true },
    @{ Ref = "E14"; Value = "  +0.38%  "; Numeric = $false },
    @{ Ref = "D15"; Value = "4.976"; Numeric = $true },
    @{ Ref = "E16"; Value = "  +0.08%  "; Numeric = $false },
    @{ Ref = "D17"; Value = "13.77"; Numeric = $true },
    @{ Ref = "E17"; Value = "  -1.46%  "; Numeric = $false },
    @{ Ref = "E18"; Value = "  -0.04%  "; Numeric = $false },
    @{ Ref = "D19"; Value = "0.000007883"; Numeric = $true },
    @{ Ref = "E19"; Value = "  -1.08%  "; Numeric = $false },
    @{ Ref = "D20"; Value = "26.502.10"; Numeric = $false },
    @{ Ref = "E20"; Value = "  -0.36%  "; Numeric = $false },
    @{ Ref = "D21"; Value = "2.080.49"; Numeric = $false },
    @{ Ref = "E21"; Value = "  -0.39%  "; Numeric = $false },
    @{ Ref = "D22"; Value = "4.581"; Numeric = $true },
    @{ Ref = "E22"; Value = "  -1.36%  "; Numeric = $false },
    @{ Ref = "E23"; Value = "  -0.31%  "; Numeric = $false },
    @{ Ref = "D24"; Value = "9.175"; Numeric = $true },
    @{ Ref = "E24"; Value = "  -3.21%  "; Numeric = $false },
    @{ Ref = "D25"; Value = "141.73"; Numeric = $true },
    @{ Ref = "E25"; Value = "  -1.41%  "; Numeric = $false },
    @{ Ref = "D26"; Value = "2.166"; Numeric = $true },
    @{ Ref = "E26"; Value = "  -0.27%  "; Numeric = $false },
    @{ Ref = "D27"; Value = "1.705"; Numeric = $true },
    @{ Ref = "D28"; Value = "16.90"; Numeric = $true },
    @{ Ref = "E28"; Value = "  -0.88%  "; Numeric = $false },
    @{ Ref = "E29"; Value = "  -1.20%  "; Numeric = $false },
    @{ Ref = "D30"; Value = "4.226"; Numeric = $true },
    @{ Ref = "E30"; Value = "  +0.81%  "; Numeric = $false },
    @{ Ref = "D31"; Value = "0.08868"; Numeric = $true },
    @{ Ref = "E31"; Value = "  +1.66%  "; Numeric = $false },
    @{ Ref = "D32"; Value = "4.021"; Numeric = $true },
    @{ Ref = "E32"; Value = "  -2.03%  "; Numeric = $false },
    @{ Ref = "D33"; Value = "0.04835"; Numeric = $true },
    @{ Ref = "E33"; Value = "  -0.81%  "; Numeric = $false },
    @{ Ref = "D34"; Value = "2.900"; Numeric = $true },
    @{ Ref = "E34"; Value = "  +0.92%  "; Numeric = $false },
    @{ Ref = "D35"; Value = "0.7228"; Numeric = $true },
    @{ Ref = "E35"; Value = "  +0.05%  "; Numeric = $false },
    @{ Ref = "E36"; Value = "  -0.37%  "; Numeric = $false },
    @{ Ref = "D37"; Value = "3.093"; Numeric = $true },
    @{ Ref = "E37"; Value = "  -0.49%  "; Numeric = $false },
    @{ Ref = "D38"; Value = "2.296"; Numeric = $true },
    @{ Ref = "E38"; Value = "  +1.03%  "; Numeric = $false },
    @{ Ref = "D39"; Value = "0.01709"; Numeric = $true },
    @{ Ref = "E39"; Value = "  -4.65%  "; Numeric = $false },
    @{ Ref = "D40"; Value = "0.4661"; Numeric = $true },
    @{ Ref = "E40"; Value = "  -4.08%  "; Numeric = $false },
    @{ Ref = "D41"; Value = "0.9020"; Numeric = $true },
    @{ Ref = "E41"; Value = "  +0.27%  "; Numeric = $false },
    @{ Ref = "D42"; Value = "106.85"; Numeric = $true },
    @{ Ref = "E42"; Value = "  -3.81%  "; Numeric = $false },
    @{ Ref = "E43"; Value = "  -1.74%  "; Numeric = $false },
    @{ Ref = "E44"; Value = "  -0.03%  "; Numeric = $false },
    @{ Ref = "D45"; Value = "7.387"; Numeric = $true },
    @{ Ref = "E45"; Value = "  -3.87%  "; Numeric = $false },
    @{ Ref = "B46"; Value = "EnergySwap"; Numeric = $false },
    @{ Ref = "C46"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; Numeric = $false },
    @{ Ref = "D46"; Value = "9.017"; Numeric = $true },
    @{ Ref = "E46"; Value = "  +0.12%  "; Numeric = $false },
    @{ Ref = "B47"; Value = "Algorand"; Numeric = $false },
    @{ Ref = "C47"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; Numeric = $false },
    @{ Ref = "D47"; Value = "0.1241"; Numeric = $true },
    @{ Ref = "E47"; Value = "  +0.42%  "; Numeric = $false },
    @{ Ref = "D48"; Value = "34.75"; Numeric = $true },
    @{ Ref = "E48"; Value = "  -1.16%  "; Numeric = $false },
    @{ Ref = "D49"; Value = "0.4038"; Numeric = $true },
    @{ Ref = "E49"; Value = "  -3.48%  "; Numeric = $false },
    @{ Ref = "B50"; Value = "EOS"; Numeric = $false },
    @{ Ref = "C50"; Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; Numeric = $false },
    @{ Ref = "D50"; Value = "0.8901"; Numeric = $true },
    @{ Ref = "E50"; Value = "  -0.22%  "; Numeric = $false },
    @{ Ref = "B51"; Value = "Cronos"; Numeric = $false },
    @{ Ref = "C51"; Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; Numeric = $false },
    @{ Ref = "D51"; Value = "0.05747"; Numeric = $true },
    @{ Ref = "E51"; Value = "  -2.27%  "; Numeric = $false }
)

foreach ($item in $updates) {
    $cell = $ws.Range($item.Ref)
    if ($item.Numeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $item.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $item.Value
    }
}
